$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shorten the question in column F (header row) - this updates the shared
# string table: the old long question is dropped and the new, shorter
# question text is appended as a new shared string.
$ws.Range("F1").Value = "How are the opportunities for development and career in the company?"

# Reset/normalize the header row formatting (D1:K1) so it uses an explicit
# font reference instead of the inherited "vertical center" alignment style.
$ws.Range("D1:K1").ClearFormats()
$ws.Range("D1:K1").Font.Name = "Calibri"

# Leave the selection on the cell that was last edited.
$ws.Range("F1").Select() | Out-Null
